# "Drop in all data files from 3.0 RMI script"
#
# The "Texas Notes" worksheet (and its averaging workings) is removed from
# the workbook. The PDiCECpDoC sheet's Perc-Decline cell (B2), which used to
# pull its value from 'Texas Notes'!B10 via a formula, is hard-coded to the
# literal value (0.13) that the older/upstream version of this input file
# used.

$wb = $excel.ActiveWorkbook

$pd = $wb.Worksheets.Item("PDiCECpDoC")

# Replace the live link to the (soon to be deleted) Texas Notes sheet with
# a literal value, and leave the selection here while this sheet is still
# the active one (so the saved sheetView reflects cell B2 as selected).
$pd.Range("B2").Value = 0.13
$pd.Range("B2").Select() | Out-Null

# Drop the Texas Notes sheet entirely - it's no longer needed now that its
# one useful output has been inlined above.
$texasNotes = $wb.Worksheets.Item("Texas Notes")
$texasNotes.Delete() | Out-Null

# Make "About" the active/selected sheet again, with the default A1
# selection, matching the target workbook state.
$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("A1").Select() | Out-Null
